$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HU Transfer")

# Seed the new "Destination" sub-header value onto J2 first so the shared string
# entry for "Destination" survives once J3's text is overwritten below.
$ws.Range("J2").Value = "Destination"

# Extend the existing grey header formatting (fill + bottom-medium border) from
# I3 across the new columns J3:Q3.
$ws.Range("I3").Copy()
$ws.Range("J3:Q3").PasteSpecial(-4122)

# Build the new row 2 "Destination" group header across J2:M2 using the same
# grey fill, then drop the bottom border since row 2 isn't the header's bottom edge.
$ws.Range("J2:M2").PasteSpecial(-4122)
$ws.Range("J2:M2").Borders.Item(9).LineStyle = -4142
$ws.Range("J2:M2").HorizontalAlignment = -4108

# Re-assert the Destination text (PasteSpecial above only copied formats) and
# fill in the rest of row 3's new headers.
$ws.Range("J2").Value = "Destination"
$ws.Range("J3").Value = "End Trough"
$ws.Range("K3").Value = "End Tray"
$ws.Range("L3").Value = "Tank"
$ws.Range("M3").Value = "Heath Unit Location"
$ws.Range("N3").Value = "Transfer Loss"
$ws.Range("O3").Value = "Final (Y/N)"
$ws.Range("P3").Value = "Crew"
$ws.Range("Q3").Value = "Comments"

# Merge the Destination sub-header across its four columns.
$ws.Range("J2:M2").Merge()

# Thin vertical borders bracket the "Destination" column group (J:M) in both
# the sub-header row and the main header row.
$ws.Range("J2").Borders.Item(7).LineStyle = 1
$ws.Range("J2").Borders.Item(7).Weight = 2
$ws.Range("J3").Borders.Item(7).LineStyle = 1
$ws.Range("J3").Borders.Item(7).Weight = 2

$ws.Range("M2").Borders.Item(10).LineStyle = 1
$ws.Range("M2").Borders.Item(10).Weight = 2
$ws.Range("M3").Borders.Item(10).LineStyle = 1
$ws.Range("M3").Borders.Item(10).Weight = 2

# Transfer Loss / Final (Y/N) sit outside the bottom-ruled header group, so
# drop their inherited bottom border while keeping the grey fill.
$ws.Range("N3:O3").Borders.Item(9).LineStyle = -4142

# Leave the new header block selected, matching the saved selection.
[void]$ws.Range("J2:Q3").Select()
